$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Ensure cells keep their original text format (prices like "1.001" or "0.1820"
# must stay as literal text, matching how the source data was authored).
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '28.909.82'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -2.91%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.881.98'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -3.39%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.002'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '329.45'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -3.73%  '
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +0.12%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4584'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -4.32%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.4090'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -1.34%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '47.85'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -1.61%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07956'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -4.01%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.9914'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -5.54%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '21.59'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -5.11%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.899.10'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -2.20%  '
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -4.55%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.052'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -5.33%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.001'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -0.05%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '88.24'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -5.03%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.06579'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -1.76%  '
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -3.98%  '
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -3.90%  '
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +0.25%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '28.893.97'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -2.86%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.399'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -4.04%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.46'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +1.51%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.196'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -3.48%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.100.65'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -3.45%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '156.29'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -3.48%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '19.52'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -3.32%  '
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -5.66%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '5.461'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -3.28%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '117.28'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -4.56%  '
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -0.84%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.09313'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.396'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -5.56%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.519'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -4.35%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.277'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -4.06%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.06045'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -3.35%  '
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -4.10%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '8.301'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -4.74%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.170'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -2.59%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.001'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +0.23%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.5776'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -5.41%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.1820'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -4.82%  '
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -7.04%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.258'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -3.02%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.07497'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +1.24%  '
$ws.Range('B47').NumberFormat = '@'
$ws.Range('B47').Value = 'Decentraland'
$ws.Range('C47').NumberFormat = '@'
$ws.Range('C47').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5441'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -5.08%  '
$ws.Range('B48').NumberFormat = '@'
$ws.Range('B48').Value = 'RenderToken'
$ws.Range('C48').NumberFormat = '@'
$ws.Range('C48').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.238'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -4.01%  '
$ws.Range('B49').NumberFormat = '@'
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').NumberFormat = '@'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '11.92'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -5.23%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.897'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -5.19%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '111.27'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -2.66%  '
